# Update "想去人数" (want-to-go count) values in column F across sheets
# to reflect the latest scraped data (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 166
$ws1.Range("F3").Value = 483
$ws1.Range("F4").Value = 16
$ws1.Range("F7").Value = 29
$ws1.Range("F8").Value = 18
$ws1.Range("F9").Value = 764

# 演出 (Performances) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 87
$ws2.Range("F3").Value = 38

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 166
$ws4.Range("F3").Value = 87
$ws4.Range("F4").Value = 483
$ws4.Range("F5").Value = 16
$ws4.Range("F8").Value = 29
$ws4.Range("F9").Value = 18
$ws4.Range("F10").Value = 764
$ws4.Range("F11").Value = 38
